$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new (blank) column before column N ("Late"/"Over Due" columns
# shift one to the right: N->O, O->P, P->Q), matching the N/O/P -> O/Q shift
# seen on the "Repayment Schedule" sheet.
$ws.Range("N1").EntireColumn.Insert()

# Make "Repayment Schedule" the active sheet/tab and select R5 on it
# (previously "Transactions" was the active/selected tab with C11 selected).
$ws.Activate()
$ws.Range("R5").Select()
